$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect so cell values can be updated.
$ws.Unprotect()

# Update the "as of" date in the confidential notice text (A59): 2021-05-06 -> 2021-05-07
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-07 for illustrative purposes only and are subject to change."

# Refresh the performance table (columns D = weight/return, E = contribution) for rows 2-56
$ws.Range("D2").Value = 0.01316704446400841
$ws.Range("E2").Value = -0.006363899334683287
$ws.Range("D3").Value = 0.01063048808224411
$ws.Range("E3").Value = 0.004975124378109541
$ws.Range("D4").Value = 0.01053842378523554
$ws.Range("E4").Value = -0.005947955390334658
$ws.Range("D5").Value = 0.01146994905993212
$ws.Range("E5").Value = 0.001176470588235334
$ws.Range("D6").Value = 0.01104249213481434
$ws.Range("E6").Value = -0.009224219489120111
$ws.Range("D7").Value = 0.0139668938529126
$ws.Range("E7").Value = -0.02406015037593989
$ws.Range("D8").Value = 0.01093127498169052
$ws.Range("E8").Value = 0.004897959183673528
$ws.Range("D9").Value = 0.01124969121460667
$ws.Range("E9").Value = -0.003250270855904569
$ws.Range("D10").Value = 0.01052144739004247
$ws.Range("E10").Value = -0.003971701625915314
$ws.Range("D11").Value = 0.01109995070316011
$ws.Range("E11").Value = 0.002647058823529447
$ws.Range("D12").Value = 0.4436062651604107
$ws.Range("E12").Value = 0.003797468354430178
$ws.Range("D13").Value = 0.01156614863269284
$ws.Range("E13").Value = 0.01151631477927051
$ws.Range("D14").Value = 0.01064180567903948
$ws.Range("E14").Value = 0.007301360057265782
$ws.Range("D15").Value = 0.01013664909899959
$ws.Range("E15").Value = 0.0004723665564476676
$ws.Range("D16").Value = 0.009863285607172746
$ws.Range("E16").Value = 0.01061388410785979
$ws.Range("D17").Value = 0.009582413325164352
$ws.Range("E17").Value = -0.006246096189881478
$ws.Range("D18").Value = 0.008239645759220307
$ws.Range("E18").Value = -0.01944106925880928
$ws.Range("D19").Value = 0.009184991560772777
$ws.Range("E19").Value = -0.01251140362309422
$ws.Range("D20").Value = 0.01274753162125663
$ws.Range("E20").Value = -0.1364179614136931
$ws.Range("D21").Value = 0.01174897135015665
$ws.Range("E21").Value = 0.03551183727909324
$ws.Range("D22").Value = 0.01164167182669277
$ws.Range("E22").Value = -0.0008319467554077642
$ws.Range("D23").Value = 0.01145906675532118
$ws.Range("E23").Value = 0.0092592592592593
$ws.Range("D24").Value = 0.01216859301595455
$ws.Range("E24").Value = 0.04283670184224664
$ws.Range("D25").Value = 0.0124928856933606
$ws.Range("E25").Value = 0.0255052264808362
$ws.Range("D26").Value = 0.01144535505151139
$ws.Range("E26").Value = 0.04603799418107135
$ws.Range("D27").Value = 0.01217425181435224
$ws.Range("E27").Value = 0.03988486842105243
$ws.Range("D28").Value = 0.01362486301899071
$ws.Range("E28").Value = -0.01317870321560355
$ws.Range("D29").Value = 0.01151456650883698
$ws.Range("E29").Value = 0.03118797845194199
$ws.Range("D30").Value = 0.007031057009129166
$ws.Range("E30").Value = -0.01106639839034218
$ws.Range("D31").Value = 0.004955257404592115
$ws.Range("E31").Value = 0.007302075326671931
$ws.Range("D32").Value = 0.009423749323936827
$ws.Range("E32").Value = 0.01336073997944487
$ws.Range("D33").Value = 0.01104553918010541
$ws.Range("E33").Value = 0.102857142857143
$ws.Range("D34").Value = 0.01011934623466819
$ws.Range("E34").Value = 0.03017561216918119
$ws.Range("D35").Value = 0.009276947034735229
$ws.Range("E35").Value = 0.01822916666666652
$ws.Range("D36").Value = 0.009574033950613925
$ws.Range("E36").Value = 0.01620859760394633
$ws.Range("D37").Value = 0.010017705509602
$ws.Range("E37").Value = 0.03490304709141268
$ws.Range("D38").Value = 0.01143784626132984
$ws.Range("E38").Value = 0.003663003663003872
$ws.Range("D39").Value = 0.01395720860180886
$ws.Range("E39").Value = 0.01896207584830356
$ws.Range("D40").Value = 0.01140509052445091
$ws.Range("E40").Value = 0.0050761421319796
$ws.Range("D41").Value = 0.01213942843959723
$ws.Range("E41").Value = 0.04460699942627655
$ws.Range("D42").Value = 0.01140095524869875
$ws.Range("E42").Value = 0.01142546245919474
$ws.Range("D43").Value = 0.01150466361164102
$ws.Range("E43").Value = 0.006715916722632498
$ws.Range("D44").Value = 0.01076760512034197
$ws.Range("E44").Value = 0.005346350534634947
$ws.Range("D45").Value = 0.01161413959602709
$ws.Range("E45").Value = -0.01077535722651679
$ws.Range("D46").Value = 0.01120920904145396
$ws.Range("E46").Value = -0.003844510892780928
$ws.Range("D47").Value = 0.0100095437811438
$ws.Range("E47").Value = 0.008447488584474971
$ws.Range("D48").Value = 0.009428428714919532
$ws.Range("E48").Value = 0.004155124653739684
$ws.Range("D49").Value = 0.009558689901112499
$ws.Range("E49").Value = 0.0369206598586016
$ws.Range("D50").Value = 0.009492416666031867
$ws.Range("E50").Value = -0.01485761452744527
$ws.Range("D51").Value = 0.008999339444110117
$ws.Range("E51").Value = 0.0356965790778383
$ws.Range("D52").Value = 0.01020542526414074
$ws.Range("E52").Value = 0.009383663894220406
$ws.Range("D53").Value = 0.008709761318412969
$ws.Range("E53").Value = 0.01274426508071369
$ws.Range("D54").Value = 0.004154319785226837
$ws.Range("E54").Value = 0.01257367387033392
$ws.Range("D55").Value = 0.004105675883615929
$ws.Range("E55").Value = -0.004028837998303558
$ws.Range("D56").Value = 0.9999999999999999
$ws.Range("E56").Value = 0.006064980417293198

# Restore sheet protection to match the original protected state.
$ws.Protect()

Write-Host "Edit complete: updated A59 notice date and 55 rows (D2:E56) of holdings data."
